$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-16 Wednesday" "2024-10-17 Thursday"

Replace-Text "748÷3=249, 1" "495÷4=123, 3"
Replace-Text "604÷3=201, 1" "735÷4=183, 3"
Replace-Text "331÷4=82, 3" "572÷3=190, 2"
Replace-Text "996÷7=142, 2" "162÷8=20, 2"
Replace-Text "379÷6=63, 1" "747÷3=249, 0"

Replace-Text "338÷4=84, 2" "433÷3=144, 1"
Replace-Text "846÷9=94, 0" "858÷6=143, 0"
Replace-Text "956÷3=318, 2" "643÷8=80, 3"
Replace-Text "540÷5=108, 0" "676÷2=338, 0"
Replace-Text "352÷6=58, 4" "783÷6=130, 3"

Replace-Text "187÷7=26, 5" "274÷6=45, 4"
Replace-Text "276÷8=34, 4" "496÷9=55, 1"
Replace-Text "798÷2=399, 0" "744÷2=372, 0"
Replace-Text "142÷3=47, 1" "644÷2=322, 0"
Replace-Text "568÷7=81, 1" "810÷6=135, 0"

Replace-Text "292÷6=48, 4" "149÷7=21, 2"
Replace-Text "786÷3=262, 0" "889÷9=98, 7"
Replace-Text "293÷2=146, 1" "348÷4=87, 0"
Replace-Text "575÷4=143, 3" "655÷6=109, 1"
Replace-Text "925÷2=462, 1" "774÷4=193, 2"

Replace-Text "655÷7=93, 4" "674÷9=74, 8"
Replace-Text "231÷4=57, 3" "803÷4=200, 3"
Replace-Text "662÷9=73, 5" "875÷9=97, 2"
Replace-Text "506÷2=253, 0" "345÷8=43, 1"
Replace-Text "170÷3=56, 2" "482÷8=60, 2"
